$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at position 347 — this shifts the existing
# row 347 (and everything below it) down by one row, preserving all of
# their original data, exactly mirroring the target diff (old row 347
# becomes new row 348, ..., old row 426 becomes new row 427).
$ws.Rows.Item(347).Insert()

# Populate the newly inserted row 347 with its data.
$ws.Range("A347").Value = 11
$ws.Range("B347").Value = "Vega Monumental Concepción"
$ws.Range("C347").Value = "Bíobío"
$ws.Range("D347").Value = 44995
$ws.Range("E347").Value = 8
$ws.Range("F347").Value = 100112017
$ws.Range("G347").Value = "Apio"
$ws.Range("H347").Value = "Americana (o)"
$ws.Range("I347").Value = "Primera"
$ws.Range("J347").Value = 150
$ws.Range("K347").Value = 7000
$ws.Range("L347").Value = 8500
$ws.Range("M347").Value = 7833
$ws.Range("N347").Value = "$/docena de matas"
$ws.Range("O347").Value = "Región de Coquimbo"
$ws.Range("P347").Value = 1306
$ws.Range("Q347").Value = 6
$ws.Range("R347").Value = "Hortaliza"
